$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F6").Value = 166
$ws.Range("F7").Value = 158
$ws.Range("F8").Value = 772
$ws.Range("F9").Value = 4156
$ws.Range("F12").Value = 170
$ws.Range("F14").Value = 5963
$ws.Range("F15").Value = 462
$ws.Range("F16").Value = 2302
$ws.Range("F19").Value = 455
$ws.Range("F20").Value = 8978
$ws.Range("F22").Value = 2235
$ws.Range("F23").Value = 194
$ws.Range("F24").Value = 2287
$ws.Range("F25").Value = 2397
$ws.Range("F26").Value = 1376
$ws.Range("F27").Value = 226
$ws.Range("F28").Value = 1937
$ws.Range("F31").Value = 324
$ws.Range("F38").Value = 1212
$ws.Range("F39").Value = 68
$ws.Range("F41").Value = 227
$ws.Range("F42").Value = 1504
$ws.Range("F43").Value = 2444
$ws.Range("F45").Value = 907
$ws.Range("F46").Value = 286
$ws.Range("F47").Value = 1245
$ws.Range("F48").Value = 5

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 167
$ws.Range("F22").Value = 32
$ws.Range("F23").Value = 32

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 685
$ws.Range("F3").Value = 882

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 685
$ws.Range("F4").Value = 882
$ws.Range("F8").Value = 166
$ws.Range("F10").Value = 167
$ws.Range("F11").Value = 158
$ws.Range("F12").Value = 772
$ws.Range("F13").Value = 4156
$ws.Range("F14").Value = 4156
$ws.Range("F16").Value = 170
$ws.Range("F19").Value = 5963
$ws.Range("F20").Value = 462
$ws.Range("F21").Value = 2302
$ws.Range("F23").Value = 455
$ws.Range("F24").Value = 8978
$ws.Range("F27").Value = 2235
$ws.Range("F28").Value = 2287
$ws.Range("F29").Value = 2397
$ws.Range("F30").Value = 1376
$ws.Range("F31").Value = 226
$ws.Range("F32").Value = 1937
$ws.Range("F35").Value = 324
$ws.Range("F40").Value = 1212
$ws.Range("F42").Value = 227
$ws.Range("F43").Value = 1504
$ws.Range("F44").Value = 2444
$ws.Range("F45").Value = 907
$ws.Range("F46").Value = 286
$ws.Range("F50").Value = 1245
$ws.Range("F51").Value = 32
